$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Localización"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "ID"
$ws.Range("E1").Value = "Tipo"

# --- Row 2 values ---
$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("B2").Value = "40.5N30.99W"
# C2 keeps its existing value/hyperlink (juan@example.com)
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "58758L"

# --- Remove now-unused columns F:I and the old E2 address value ---
$ws.Range("F1:I2").ClearContents()
$ws.Range("E2").ClearContents()

# --- Selection ends up on E3 ---
[void]$ws.Range("E3").Select()
